# The "growth process" was re-run, producing new voxel_total (B) and
# voxel_height (C) values, and a uniform voxel_depth (D) of 7 for every
# space row (2-20). Update the data table accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# row -> (voxel_total, voxel_height, voxel_depth)
$data = @{
    2  = @(408, 1, 7)
    3  = @(360, 1, 7)
    4  = @(1530, 1, 7)
    5  = @(90, 1, 7)
    6  = @(102, 1, 7)
    7  = @(26, 1, 7)
    8  = @(24, 1, 7)
    9  = @(64, 1, 7)
    10 = @(36, 1, 7)
    11 = @(114, 1, 7)
    12 = @(52, 1, 7)
    13 = @(256, 1, 7)
    14 = @(64, 1, 7)
    15 = @(76, 1, 7)
    16 = @(2, 1, 7)
    17 = @(178, 1, 7)
    18 = @(4, 1, 7)
    19 = @(128, 1, 7)
    20 = @(52, 1, 7)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
}

# Update the active selection to B20, matching the new cursor position
$ws.Range("B20").Select()
